$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Collapse the "Best weeks: ..." paragraph (currently split across many
#    runs with a _GoBack bookmark in the middle) into a single run.
# ---------------------------------------------------------------------------
$bw = $d.Paragraphs(8).Range
$bwBody = $d.Range($bw.Start, $bw.End - 1)
$bwBody.Delete()
$d.Range($bw.Start, $bw.Start).InsertAfter("Best weeks: July 22, August 19, September 9, September 16")

# ---------------------------------------------------------------------------
# 2) Append the new "Interested Parties" entries at the end of the document,
#    after the "POC: Charles Baldner: baldner@stanford.edu " paragraph.
# ---------------------------------------------------------------------------

function Add-Para($level) {
    $d.Paragraphs.Last.Range.InsertParagraphAfter()
    $d.Paragraphs.Last.Range.ListFormat.ListLevelNumber = $level
}

function Add-PocHyperlink($prefixText, $email) {
    $d.Paragraphs.Last.Range.InsertAfter($prefixText)
    $prefixEnd = $d.Paragraphs.Last.Range.End
    $d.Paragraphs.Last.Range.InsertAfter($email)
    $emailRange = $d.Range($prefixEnd - 1, $d.Paragraphs.Last.Range.End)
    $d.Hyperlinks.Add($emailRange, "mailto:" + $email, $null, $null, $email) | Out-Null
}

# Laura Kay Schaefer
Add-Para 2
$d.Paragraphs.Last.Range.InsertAfter("Laura Kay Schaefer")

Add-Para 3
$d.Paragraphs.Last.Range.InsertAfter("Planetary science, seismology")

Add-Para 3
Add-PocHyperlink "POC: Laura Schaefer: " "lkschaef@stanford.edu"

# Sigrid Close
Add-Para 2
$d.Paragraphs.Last.Range.InsertAfter("Sigrid Close")

Add-Para 3
$d.Paragraphs.Last.Range.InsertAfter("CubeSats, micrometeroites, plasma")

Add-Para 3
Add-PocHyperlink "POC: Nicholas Lee: " "nnlee@stanford.edu"
$d.Paragraphs.Last.Range.InsertAfter("`t")

# Monica Bobra
Add-Para 2
$d.Paragraphs.Last.Range.InsertAfter("Monica Bobra")

Add-Para 3
$d.Paragraphs.Last.Range.InsertAfter("Heliophysics and sun weather")

Add-Para 3
Add-PocHyperlink "POC: Monica Bobra: " "mbobra@stanford.edu"
$d.Paragraphs.Last.Range.InsertAfter("`t")

# Alexandra Koenig
Add-Para 1
$d.Paragraphs.Last.Range.InsertAfter("Alexandra Koenig Remote sensing")
$d.Paragraphs.Last.Range.InsertAfter("X")
$p = $d.Paragraphs.Last
$pEnd = $p.Range.End
$markRange = $d.Range($pEnd - 2, $pEnd - 1)
$d.Bookmarks.Add("_GoBack", $markRange) | Out-Null
$delRange = $d.Range($pEnd - 2, $pEnd - 1)
$delRange.Delete()
